$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- J1 header: copy the bold/border/center-top style used by the other
#     header cells (B1:I1, style index 1) onto the new header cell, then set its text
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("J1").Value = "epl"

# --- J2:J49 data: format range as Text first so values such as "0.490" or
#     "inf" are stored verbatim as strings rather than being parsed as numbers
#     (this matches the source workbook, which keeps these as literal strings)
$ws.Range("J2:J49").NumberFormat = "@"
$ws.Range("J2").Value = "inf"
$ws.Range("J3").Value = "0.031"
$ws.Range("J4").Value = "0.024"
$ws.Range("J5").Value = "0.067"
$ws.Range("J6").Value = "0.129"
$ws.Range("J7").Value = "0.288"
$ws.Range("J8").Value = "inf"
$ws.Range("J9").Value = "0.042"
$ws.Range("J10").Value = "0.046"
$ws.Range("J11").Value = "0.085"
$ws.Range("J12").Value = "0.168"
$ws.Range("J13").Value = "0.393"
$ws.Range("J14").Value = "inf"
$ws.Range("J15").Value = "0.047"
$ws.Range("J16").Value = "0.066"
$ws.Range("J17").Value = "0.107"
$ws.Range("J18").Value = "0.153"
$ws.Range("J19").Value = "0.490"
$ws.Range("J20").Value = "inf"
$ws.Range("J21").Value = "0.081"
$ws.Range("J22").Value = "0.085"
$ws.Range("J23").Value = "0.141"
$ws.Range("J24").Value = "0.166"
$ws.Range("J25").Value = "0.488"
$ws.Range("J26").Value = "inf"
$ws.Range("J27").Value = "0.130"
$ws.Range("J28").Value = "0.157"
$ws.Range("J29").Value = "0.270"
$ws.Range("J30").Value = "0.262"
$ws.Range("J31").Value = "0.715"
$ws.Range("J32").Value = "inf"
$ws.Range("J33").Value = "0.275"
$ws.Range("J34").Value = "0.304"
$ws.Range("J35").Value = "0.531"
$ws.Range("J36").Value = "0.460"
$ws.Range("J37").Value = "0.450"
$ws.Range("J38").Value = "inf"
$ws.Range("J39").Value = "0.426"
$ws.Range("J40").Value = "0.427"
$ws.Range("J41").Value = "0.716"
$ws.Range("J42").Value = "0.659"
$ws.Range("J43").Value = "0.662"
$ws.Range("J44").Value = "inf"
$ws.Range("J45").Value = "0.758"
$ws.Range("J46").Value = "0.524"
$ws.Range("J47").Value = "0.818"
$ws.Range("J48").Value = "0.819"
$ws.Range("J49").Value = "0.850"
